$wb = $excel.ActiveWorkbook

# --- Update the "Date" metadata value ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-07-25T07:22:51+00:00"

# --- Update the canonical terminology URLs on the Elements sheet (column Z) ---
$elements = $wb.Worksheets.Item("Elements")

$elements.Range("Z5").Value  = "https://mos.esante.gouv.fr/NOS/TRE_R23-ModeExercice/FHIR/TRE-R23-ModeExercice?vs"
$elements.Range("Z9").Value  = "https://mos.esante.gouv.fr/NOS/TRE_R22-GenreActivite/FHIR/TRE-R22-GenreActivite?vs"
$elements.Range("Z10").Value = "https://mos.esante.gouv.fr/NOS/TRE_R25-MotifFinActivite/FHIR/TRE-R25-MotifFinActivite?vs"
$elements.Range("Z11").Value = "https://mos.esante.gouv.fr/NOS/TRE_R32-StatutHospitalier/FHIR/TRE-R32-StatutHospitalier?vs"
$elements.Range("Z12").Value = "https://mos.esante.gouv.fr/NOS/TRE_R06-SectionTableauCNOP/FHIR/TRE-R06-SectionTableauCNOP?vs"
$elements.Range("Z13").Value = "https://mos.esante.gouv.fr/NOS/TRE_G05-SousSectionTableauCNOP/FHIR/TRE-G05-SousSectionTableauCNOP?vs"
$elements.Range("Z14").Value = "https://mos.esante.gouv.fr/NOS/TRE_R24-TypeActiviteLiberale/FHIR/TRE-R24-TypeActiviteLiberale?vs"
$elements.Range("Z15").Value = "https://mos.esante.gouv.fr/NOS/TRE_R34-StatutProfessionnelSSA/FHIR/TRE-R34-StatutProfessionnelSSA?vs"

# --- Resize column Z to fit the new (longer) content, as in the original export ---
# (ColumnWidth is quantized to the renderer's pixel grid on save; 88.3 is the input
# that lands on the grid point closest to the authored width of 89.234375 chars.)
$elements.Columns.Item(26).ColumnWidth = 88.3
